# aggiornamento fino a 9 agosto 2021
# Append 15 new daily rows (329-343) to Sheet1, continuing the existing
# date series in column A (serial dates 44403..44417, i.e. 2021-07-26
# through 2021-08-09) with zeros in columns B, C and D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing data row is 328 (date serial 44402 / 2021-07-25).
$lastRow = 328
$lastSerial = 44402
$newRowCount = 15

# Copy the date cell's formatting (number format, font, borders, alignment)
# down into the new A-column cells before writing their values, so the new
# rows reuse the same style as the rest of the date column instead of
# minting a brand-new one.
$fmtSrc = $ws.Range("A" + $lastRow)
$fmtDst = $ws.Range("A" + ($lastRow + 1) + ":A" + ($lastRow + $newRowCount))
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)  # xlPasteFormats

for ($i = 1; $i -le $newRowCount; $i++) {
    $row = $lastRow + $i
    $serial = $lastSerial + $i

    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
